$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds figures that are stored as TEXT in the
# workbook (e.g. "14.00", "583.78"), not numbers. Force the Text number
# format before writing so Excel does not auto-convert numeric-looking
# strings into real numbers and silently drop trailing zeros / precision.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.322.10"
$ws.Range("E2").Value = "  -2.63%  "

$ws.Range("D3").Value = "3.007.28"
$ws.Range("E3").Value = "  -2.47%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "583.78"
$ws.Range("E5").Value = "  -1.80%  "

$ws.Range("D6").Value = "146.31"
$ws.Range("E6").Value = "  -6.12%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").Value = "0.525"
$ws.Range("E8").Value = "  -2.65%  "

$ws.Range("D9").Value = "3.005.83"
$ws.Range("E9").Value = "  -2.42%  "

$ws.Range("D10").Value = "0.149"
$ws.Range("E10").Value = "  -5.41%  "

$ws.Range("D11").Value = "5.82"
$ws.Range("E11").Value = "  -1.67%  "

$ws.Range("D12").Value = "0.463"
$ws.Range("E12").Value = "  +2.29%  "

$ws.Range("D13").Value = "0.0000229"
$ws.Range("E13").Value = "  -4.62%  "

$ws.Range("D14").Value = "34.51"
$ws.Range("E14").Value = "  -6.89%  "

$ws.Range("D15").Value = "0.122"
$ws.Range("E15").Value = "  +1.57%  "

$ws.Range("D16").Value = "3.506.66"
$ws.Range("E16").Value = "  -2.48%  "

$ws.Range("D17").Value = "7.12"
$ws.Range("E17").Value = "  -1.46%  "

$ws.Range("D18").Value = "62.352.62"
$ws.Range("E18").Value = "  -2.56%  "

$ws.Range("D19").Value = "3.013.10"
$ws.Range("E19").Value = "  -2.47%  "

$ws.Range("D20").Value = "461.19"
$ws.Range("E20").Value = "  -5.30%  "

$ws.Range("D21").Value = "14.00"
$ws.Range("E21").Value = "  -4.30%  "

$ws.Range("D22").Value = "0.687"
$ws.Range("E22").Value = "  -3.66%  "

$ws.Range("D23").Value = "7.45"
$ws.Range("E23").Value = "  -2.37%  "

$ws.Range("D24").Value = "81.46"
$ws.Range("E24").Value = "  -0.80%  "

$ws.Range("D25").Value = "2.23"
$ws.Range("E25").Value = "  -8.76%  "

$ws.Range("D26").Value = "12.26"
$ws.Range("E26").Value = "  -5.38%  "

$ws.Range("D27").Value = "10.14"
$ws.Range("E27").Value = "  -4.69%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").Value = "2.63"
$ws.Range("E30").Value = "  -2.97%  "

$ws.Range("D31").Value = "7.07"
$ws.Range("E31").Value = "  -6.91%  "

$ws.Range("D32").Value = "2.10"
$ws.Range("E32").Value = "  -7.14%  "

$ws.Range("D33").Value = "28.40"
$ws.Range("E33").Value = "  +3.64%  "

$ws.Range("D34").Value = "0.108"
$ws.Range("E34").Value = "  -3.89%  "

$ws.Range("D35").Value = "0.0₃0796"
$ws.Range("E35").Value = "  -4.11%  "

$ws.Range("E36").Value = "  -3.96%  "

$ws.Range("D37").Value = "5.76"
$ws.Range("E37").Value = "  -5.51%  "

$ws.Range("D38").Value = "2.11"
$ws.Range("E38").Value = "  -6.18%  "

$ws.Range("D39").Value = "50.38"
$ws.Range("E39").Value = "  -0.76%  "

$ws.Range("D40").Value = "9.17"
$ws.Range("E40").Value = "  -1.57%  "

$ws.Range("D41").Value = "2.87"
$ws.Range("E41").Value = "  -12.02%  "

$ws.Range("E42").Value = "  +1.15%  "

$ws.Range("D43").Value = "395.64"
$ws.Range("E43").Value = "  -10.49%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0358"
$ws.Range("E44").Value = "  -2.78%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.273"
$ws.Range("E45").Value = "  -6.54%  "

$ws.Range("D46").Value = "2.736.50"
$ws.Range("E46").Value = "  -3.81%  "

$ws.Range("D47").Value = "36.79"
$ws.Range("E47").Value = "  -7.63%  "

$ws.Range("D48").Value = "129.26"
$ws.Range("E48").Value = "  -2.75%  "

$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("D50").Value = "0.109"
$ws.Range("E50").Value = "  -0.89%  "

$ws.Range("D51").Value = "2.19"
$ws.Range("E51").Value = "  -2.98%  "

